$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16504777985325565"
$ws1.Range("B2").Value = "go_stims-16504777984965587.csv"
$ws1.Range("B3").Value = "GNG_stims-16504777985165966.csv"
$ws1.Range("B4").Value = "go_stims-16504777985175536.csv"
$ws1.Range("B5").Value = "GNG_stims-16504777985315874.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16504778005236228"
$ws2.Range("B2").Value = "ZB-match_5-1650477799088553.csv"
$ws2.Range("B3").Value = "ZB-match_3-16504777990375855.csv"
$ws2.Range("B4").Value = "TB-1650477800107552.csv"
$ws2.Range("B5").Value = "OB-16504777995305877.csv"
$ws2.Range("B6").Value = "OB-16504777993665879.csv"
$ws2.Range("B7").Value = "ZB-match_8-1650477798954584.csv"
$ws2.Range("B8").Value = "OB-16504777994455872.csv"
$ws2.Range("B9").Value = "TB-1650477800511588.csv"
$ws2.Range("B10").Value = "TB-1650477800235585.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-1650477800524555"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16504778005715551"
$ws4.Range("B2").Value = "MM_stims-16504778005395849.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778005255568.csv"
$ws4.Range("B4").Value = "MM_stims-16504778005555844.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778005395849.csv"
$ws4.Range("B6").Value = "MM_stims-16504778005715551.csv"
$ws4.Range("B7").Value = "ZM_stims-1650477800556553.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16504778006355858"
$ws5.Range("B2").Value = "SAT_stims-16504778005745575.csv"
$ws5.Range("B3").Value = "SAT_stims-1650477800587585.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504778006035852.csv"
$ws5.Range("B5").Value = "vSAT_stims-1650477800619553.csv"
